$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update formulas for G2 and H2 (swap/shift values to reflect new LDO 3.3v ref voltages)
$ws.Range("G2").Formula = "=`$B`$2+1408+384"
$ws.Range("H2").Formula = "=`$B`$2+512"

# Update linker heap formula in C3
$ws.Range("C3").Formula = "=16*1024"

# Recalculate the workbook so dependent cells' cached values update
$excel.CalculateFullRebuild()

# Update selection to reflect new active cell G2
$ws.Range("G2").Select()
